$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Update the Status column ("N/A" -> "IP") for the rows that are now "In Progress":
# POWER (row 8), LN (row 20), LG (row 21), LOG (row 22), LOGN (row 23)
$ws.Range("D8").Value = "IP"
$ws.Range("D20").Value = "IP"
$ws.Range("D21").Value = "IP"
$ws.Range("D22").Value = "IP"
$ws.Range("D23").Value = "IP"

# Widen column B to fit the longer entries
$ws.Columns.Item(2).ColumnWidth = 25.736979166666668

# Move the active selection
$ws.Range("H8").Select()

$wb.Save()
